$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8429523706436157
$ws.Range("B1").Value = 1.84393298625946
$ws.Range("D1").Value = 1.938547730445862
$ws.Range("E1").Value = 1.144961476325989
